$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 810, shifting the existing
# rows 810-881 down to 812-883.
$ws.Rows("810:811").Insert()

# New row 810: Cebolla, 1a (guarda)
$ws.Range("A810").Value = 11
$ws.Range("B810").Value = "Vega Monumental Concepción"
$ws.Range("C810").Value = "Bíobío"
$ws.Range("D810").Value = 45132
$ws.Range("E810").Value = 8
$ws.Range("F810").Value = 100112004
$ws.Range("G810").Value = "Cebolla"
$ws.Range("H810").Value = "Sin especificar"
$ws.Range("I810").Value = "1a (guarda)"
$ws.Range("J810").Value = 400
$ws.Range("K810").Value = 9000
$ws.Range("L810").Value = 9500
$ws.Range("M810").Value = 9250
$ws.Range("N810").Value = "`$/malla 18 kilos"
$ws.Range("O810").Value = "Región de O'Higgins"
$ws.Range("P810").Value = 514
$ws.Range("Q810").Value = 18
$ws.Range("R810").Value = "Hortaliza"

# New row 811: Cebolla, 2a (guarda)
$ws.Range("A811").Value = 11
$ws.Range("B811").Value = "Vega Monumental Concepción"
$ws.Range("C811").Value = "Bíobío"
$ws.Range("D811").Value = 45132
$ws.Range("E811").Value = 8
$ws.Range("F811").Value = 100112004
$ws.Range("G811").Value = "Cebolla"
$ws.Range("H811").Value = "Sin especificar"
$ws.Range("I811").Value = "2a (guarda)"
$ws.Range("J811").Value = 200
$ws.Range("K811").Value = 8000
$ws.Range("L811").Value = 8000
$ws.Range("M811").Value = 8000
$ws.Range("N811").Value = "`$/malla 18 kilos"
$ws.Range("O811").Value = "Región de O'Higgins"
$ws.Range("P811").Value = 444
$ws.Range("Q811").Value = 18
$ws.Range("R811").Value = "Hortaliza"

# Keep the date columns formatted the same way as the rest of column D.
$ws.Range("D810:D811").NumberFormat = $ws.Range("D812").NumberFormat
